$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRange, [string]$value)
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $value
    $cellRange.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "303.50"
Set-TextValue $ws.Range("E2") "-4.74%"
Set-TextValue $ws.Range("D3") "35.21"
Set-TextValue $ws.Range("E3") "-2.50%"
Set-TextValue $ws.Range("D4") "5.064"
Set-TextValue $ws.Range("E4") "-1.98%"
Set-TextValue $ws.Range("D5") "0.08003"
Set-TextValue $ws.Range("D6") "1.927"
Set-TextValue $ws.Range("E6") "-10.29%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D7") "4.049"
Set-TextValue $ws.Range("E7") "-2.14%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue $ws.Range("D8") "7.749"
Set-TextValue $ws.Range("E8") "-3.26%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D9") "2.989"
Set-TextValue $ws.Range("E9") "6.71%"
$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D10") "0.9213"
Set-TextValue $ws.Range("E10") "-0.64%"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D11") "0.1218"
Set-TextValue $ws.Range("E11") "19.48%"
$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D12") "0.1847"
Set-TextValue $ws.Range("E12") "-2.38%"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D13") "0.09624"
Set-TextValue $ws.Range("E13") "4.58%"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D14") "0.03589"
Set-TextValue $ws.Range("E14") "-0.92%"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D15") "0.09859"
Set-TextValue $ws.Range("E15") "-0.66%"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D16") "0.001390"
Set-TextValue $ws.Range("E16") "-3.88%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D17") "0.005746"
Set-TextValue $ws.Range("E17") "0.66%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D18") "3.500"
Set-TextValue $ws.Range("E18") "1.07%"
Set-TextValue $ws.Range("D19") "0.3412"
Set-TextValue $ws.Range("E19") "1.10%"
Set-TextValue $ws.Range("D20") "0.1283"
Set-TextValue $ws.Range("E20") "-1.46%"
Set-TextValue $ws.Range("D21") "5.041"
Set-TextValue $ws.Range("E21") "-3.19%"
Set-TextValue $ws.Range("D22") "0.2466"
Set-TextValue $ws.Range("E22") "12.51%"
Set-TextValue $ws.Range("D23") "0.04501"
Set-TextValue $ws.Range("E23") "-2.18%"
Set-TextValue $ws.Range("D24") "0.001214"
Set-TextValue $ws.Range("E24") "-2.75%"
Set-TextValue $ws.Range("E25") "2.54%"
Set-TextValue $ws.Range("D26") "0.0001250"
Set-TextValue $ws.Range("E26") "-0.10%"
Set-TextValue $ws.Range("E27") "-6.93%"
Set-TextValue $ws.Range("D39") "0.01930"
Set-TextValue $ws.Range("E39") "-3.97%"
Set-TextValue $ws.Range("D40") "0.04751"
Set-TextValue $ws.Range("E40") "-3.55%"
Set-TextValue $ws.Range("D41") "0.007547"
Set-TextValue $ws.Range("E41") "-3.07%"
Set-TextValue $ws.Range("D42") "0.009557"
Set-TextValue $ws.Range("E42") "22.28%"
Set-TextValue $ws.Range("D43") "0.1332"
Set-TextValue $ws.Range("E43") "-4.95%"
Set-TextValue $ws.Range("E44") "0.21%"
Set-TextValue $ws.Range("D45") "0.01116"
Set-TextValue $ws.Range("E45") "-6.55%"
Set-TextValue $ws.Range("D46") "0.00006275"
Set-TextValue $ws.Range("E46") "-2.88%"
Set-TextValue $ws.Range("E47") "-0.09%"
Set-TextValue $ws.Range("E48") "57.92%"
Set-TextValue $ws.Range("E49") "-31.43%"
Set-TextValue $ws.Range("D50") "0.00002101"
Set-TextValue $ws.Range("E50") "-0.09%"
Set-TextValue $ws.Range("D51") "0.0002000"
Set-TextValue $ws.Range("E51") "-0.09%"
